$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 51; existing rows 51-70 shift down to 52-71.
$ws.Rows("51:51").Insert()

# Populate the newly inserted row 51 with the new price-report record.
$ws.Range("A51").Value = 8
$ws.Range("B51").Value = "Terminal La Palmera de La Serena"
$ws.Range("C51").Value = "Coquimbo"
$ws.Range("D51").Value = 44588
$ws.Range("E51").Value = 4
$ws.Range("F51").Value = "Fruta"
$ws.Range("G51").Value = 100109
$ws.Range("H51").Value = "Uva"
$ws.Range("I51").Value = 100109001
$ws.Range("J51").Value = "Uva"
$ws.Range("K51").Value = "Flame Seedless"
$ws.Range("L51").Value = "Primera"
$ws.Range("M51").Value = 500
$ws.Range("N51").Value = 7500
$ws.Range("O51").Value = 8000
$ws.Range("P51").Value = 7750
$ws.Range("Q51").Value = "$/caja 15 kilos"
$ws.Range("R51").Value = "Provincia del Elquí"
$ws.Range("S51").Value = 517
$ws.Range("T51").Value = 15
